$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 10 de Septiembre de 2020 a las 00:44'
$ws.Cells.Item(4, 2).Value = 6546726
$ws.Cells.Item(4, 3).Value = 32495
$ws.Cells.Item(4, 4).Value = 3838941
$ws.Cells.Item(4, 5).Value = 2512637
$ws.Cells.Item(4, 7).Value = 1118
$ws.Cells.Item(4, 8).Value = 195148
$ws.Cells.Item(6, 4).Value = 3453336
$ws.Cells.Item(6, 5).Value = 616014
$ws.Cells.Item(8, 2).Value = 702776
$ws.Cells.Item(8, 3).Value = 6586
$ws.Cells.Item(8, 4).Value = 536959
$ws.Cells.Item(8, 5).Value = 135581
$ws.Cells.Item(8, 7).Value = 113
$ws.Cells.Item(8, 8).Value = 30236
$ws.Cells.Item(9, 2).Value = 686856
$ws.Cells.Item(9, 3).Value = 7343
$ws.Cells.Item(9, 4).Value = 552885
$ws.Cells.Item(9, 5).Value = 111918
$ws.Cells.Item(9, 7).Value = 236
$ws.Cells.Item(9, 8).Value = 22053
$ws.Cells.Item(23, 2).Value = 273821
$ws.Cells.Item(23, 3).Value = 4243
$ws.Cells.Item(23, 4).Value = 209993
$ws.Cells.Item(23, 5).Value = 56096
$ws.Cells.Item(23, 7).Value = 75
$ws.Cells.Item(23, 8).Value = 7732
$ws.Cells.Item(24, 2).Value = 256349
$ws.Cells.Item(24, 3).Value = 1393
$ws.Cells.Item(24, 4).Value = 231900
$ws.Cells.Item(24, 5).Value = 15039
$ws.Cells.Item(24, 7).Value = 1
$ws.Cells.Item(24, 8).Value = 9410
$ws.Cells.Item(29, 2).Value = 134194
$ws.Cells.Item(29, 3).Value = 446
$ws.Cells.Item(29, 4).Value = 118149
$ws.Cells.Item(29, 5).Value = 6890
$ws.Cells.Item(29, 7).Value = 2
$ws.Cells.Item(29, 8).Value = 9155
$ws.Cells.Item(35, 2).Value = 100403
$ws.Cells.Item(35, 3).Value = 175
$ws.Cells.Item(35, 4).Value = 80689
$ws.Cells.Item(35, 5).Value = 14137
$ws.Cells.Item(35, 7).Value = 17
$ws.Cells.Item(35, 8).Value = 5577
$ws.Cells.Item(48, 2).Value = 72726
$ws.Cells.Item(48, 3).Value = 492
$ws.Cells.Item(48, 4).Value = 64100
$ws.Cells.Item(48, 5).Value = 7233
$ws.Cells.Item(48, 7).Value = 16
$ws.Cells.Item(48, 8).Value = 1393
$ws.Cells.Item(55, 1).Value = 'Nigeria'
$ws.Cells.Item(55, 2).Value = 55632
$ws.Cells.Item(55, 3).Value = 176
$ws.Cells.Item(55, 4).Value = 43610
$ws.Cells.Item(55, 5).Value = 10952
$ws.Cells.Item(55, 7).Value = 3
$ws.Cells.Item(55, 8).Value = 1070
$ws.Cells.Item(56, 1).Value = 'Venezuela'
$ws.Cells.Item(56, 2).Value = 55563
$ws.Cells.Item(56, 4).Value = 44435
$ws.Cells.Item(56, 5).Value = 10684
$ws.Cells.Item(56, 8).Value = 444
$ws.Cells.Item(60, 1).Value = 'Ghana'
$ws.Cells.Item(60, 2).Value = 45313
$ws.Cells.Item(60, 3).Value = 301
$ws.Cells.Item(60, 4).Value = 44188
$ws.Cells.Item(60, 5).Value = 842
$ws.Cells.Item(60, 7).Value = 0
$ws.Cells.Item(60, 8).Value = 283
$ws.Cells.Item(61, 1).Value = 'Suiza'
$ws.Cells.Item(61, 2).Value = 45306
$ws.Cells.Item(61, 3).Value = 469
$ws.Cells.Item(61, 4).Value = 38100
$ws.Cells.Item(61, 5).Value = 5187
$ws.Cells.Item(61, 7).Value = 1
$ws.Cells.Item(61, 8).Value = 2019
$ws.Cells.Item(63, 2).Value = 44930
$ws.Cells.Item(63, 3).Value = 649
$ws.Cells.Item(63, 4).Value = 42212
$ws.Cells.Item(63, 5).Value = 2352
$ws.Cells.Item(63, 7).Value = 8
$ws.Cells.Item(63, 8).Value = 366
$ws.Cells.Item(84, 2).Value = 17435
$ws.Cells.Item(84, 3).Value = 122
$ws.Cells.Item(84, 4).Value = 12474
$ws.Cells.Item(84, 5).Value = 4259
$ws.Cells.Item(84, 7).Value = 10
$ws.Cells.Item(84, 8).Value = 702
$ws.Cells.Item(92, 2).Value = 11746
$ws.Cells.Item(92, 3).Value = 122
$ws.Cells.Item(92, 5).Value = 2134
$ws.Cells.Item(105, 2).Value = 7429
$ws.Cells.Item(105, 3).Value = 41
$ws.Cells.Item(105, 4).Value = 5542
$ws.Cells.Item(105, 5).Value = 1665
$ws.Cells.Item(105, 7).Value = 4
$ws.Cells.Item(105, 8).Value = 222
$ws.Cells.Item(106, 2).Value = 7191
$ws.Cells.Item(106, 3).Value = 26
$ws.Cells.Item(106, 4).Value = 6701
$ws.Cells.Item(106, 5).Value = 329
$ws.Cells.Item(106, 7).Value = 1
$ws.Cells.Item(106, 8).Value = 161
$ws.Cells.Item(112, 2).Value = 4990
$ws.Cells.Item(112, 3).Value = 5
$ws.Cells.Item(112, 4).Value = 4479
$ws.Cells.Item(112, 5).Value = 428
$ws.Cells.Item(136, 1).Value = 'Aruba'
$ws.Cells.Item(136, 2).Value = 2730
$ws.Cells.Item(136, 3).Value = 141
$ws.Cells.Item(136, 4).Value = 1351
$ws.Cells.Item(136, 5).Value = 1364
$ws.Cells.Item(136, 8).Value = 15
$ws.Cells.Item(137, 1).Value = 'Jordania'
$ws.Cells.Item(137, 2).Value = 2659
$ws.Cells.Item(137, 3).Value = 78
$ws.Cells.Item(137, 4).Value = 1919
$ws.Cells.Item(137, 5).Value = 721
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 19
$ws.Cells.Item(138, 1).Value = 'Bahamas'
$ws.Cells.Item(138, 2).Value = 2657
$ws.Cells.Item(138, 3).Value = 72
$ws.Cells.Item(138, 4).Value = 1088
$ws.Cells.Item(138, 5).Value = 1506
$ws.Cells.Item(138, 7).Value = 4
$ws.Cells.Item(138, 8).Value = 63
$ws.Cells.Item(166, 2).Value = 1048
$ws.Cells.Item(166, 3).Value = 3
$ws.Cells.Item(166, 4).Value = 931
$ws.Cells.Item(166, 5).Value = 38
$ws.Cells.Item(169, 2).Value = 795
$ws.Cells.Item(169, 3).Value = 22
$ws.Cells.Item(169, 5).Value = 267
$ws.Cells.Item(214, 1).Value = 'Montserrat'
$ws.Cells.Item(214, 4).Value = 12
$ws.Cells.Item(214, 8).Value = 1
$ws.Cells.Item(215, 1).Value = 'Islas Malvinas'
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 8).Value = 0
